$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 160.9723563333334
$ws.Range("H2").Value = 482.917069
$ws.Range("I2").Value = 0.3931645655589854
$ws.Range("J2").Value = 0.3931645655589854
$ws.Range("M2").Value = 77.08952333333333
$ws.Range("N2").Value = 231.26857
$ws.Range("O2").Value = 0.2403816673726824
$ws.Range("P2").Value = 0.2403816673726824
$ws.Range("Q2").Value = 12409.28221958015
$ws.Range("R2").Value = 111683.5399762213
$ws.Range("S2").Value = 0.09450955382092521
$ws.Range("T2").Value = 0.09450955382092521

$ws.Range("G3").Value = 160.9723563333334
$ws.Range("H3").Value = 482.917069
$ws.Range("I3").Value = 0.3931645655589854
$ws.Range("J3").Value = 0.3931645655589854
$ws.Range("M3").Value = 101.5800373333333
$ws.Range("N3").Value = 304.740112
$ws.Range("O3").Value = 0.3167483425780597
$ws.Range("P3").Value = 0.3167483425780597
$ws.Range("Q3").Value = 16351.57796597464
$ws.Range("R3").Value = 147164.2016937717
$ws.Range("S3").Value = 0.1245342245012315
$ws.Range("T3").Value = 0.1245342245012315

$ws.Range("G4").Value = 160.9723563333334
$ws.Range("H4").Value = 482.917069
$ws.Range("I4").Value = 0.3931645655589854
$ws.Range("J4").Value = 0.3931645655589854
$ws.Range("M4").Value = 142.0267893333333
$ws.Range("N4").Value = 426.080368
$ws.Range("O4").Value = 0.4428699900492579
$ws.Range("P4").Value = 0.4428699900492579
$ws.Range("Q4").Value = 22862.3869414446
$ws.Range("R4").Value = 205761.4824730014
$ws.Range("S4").Value = 0.1741207872368287
$ws.Range("T4").Value = 0.1741207872368287

$ws.Range("G5").Value = 89.97721833333333
$ws.Range("H5").Value = 269.931655
$ws.Range("I5").Value = 0.2197635343237224
$ws.Range("J5").Value = 0.2197635343237224
$ws.Range("M5").Value = 77.08952333333333
$ws.Range("N5").Value = 231.26857
$ws.Range("O5").Value = 0.2403816673726824
$ws.Range("P5").Value = 0.2403816673726824
$ws.Range("Q5").Value = 6936.300872175927
$ws.Range("R5").Value = 62426.70784958335
$ws.Range("S5").Value = 0.0528271248084501
$ws.Range("T5").Value = 0.05282712480845011

$ws.Range("G6").Value = 89.97721833333333
$ws.Range("H6").Value = 269.931655
$ws.Range("I6").Value = 0.2197635343237224
$ws.Range("J6").Value = 0.2197635343237224
$ws.Range("M6").Value = 101.5800373333333
$ws.Range("N6").Value = 304.740112
$ws.Range("O6").Value = 0.3167483425780597
$ws.Range("P6").Value = 0.3167483425780597
$ws.Range("Q6").Value = 9139.889197449484
$ws.Range("R6").Value = 82259.00277704536
$ws.Range("S6").Value = 0.0696097352561356
$ws.Range("T6").Value = 0.06960973525613559

$ws.Range("G7").Value = 89.97721833333333
$ws.Range("H7").Value = 269.931655
$ws.Range("I7").Value = 0.2197635343237224
$ws.Range("J7").Value = 0.2197635343237224
$ws.Range("M7").Value = 142.0267893333333
$ws.Range("N7").Value = 426.080368
$ws.Range("O7").Value = 0.4428699900492579
$ws.Range("P7").Value = 0.4428699900492579
$ws.Range("Q7").Value = 12779.17543302767
$ws.Range("R7").Value = 115012.578897249
$ws.Range("S7").Value = 0.09732667425913667
$ws.Range("T7").Value = 0.09732667425913669

$ws.Range("G8").Value = 158.477852
$ws.Range("H8").Value = 475.433556
$ws.Range("I8").Value = 0.3870719001172923
$ws.Range("J8").Value = 0.3870719001172923
$ws.Range("M8").Value = 77.08952333333333
$ws.Range("N8").Value = 231.26857
$ws.Range("O8").Value = 0.2403816673726824
$ws.Range("P8").Value = 0.2403816673726824
$ws.Range("Q8").Value = 12216.98206957055
$ws.Range("R8").Value = 109952.8386261349
$ws.Range("S8").Value = 0.09304498874330709
$ws.Range("T8").Value = 0.0930449887433071

$ws.Range("G9").Value = 158.477852
$ws.Range("H9").Value = 475.433556
$ws.Range("I9").Value = 0.3870719001172923
$ws.Range("J9").Value = 0.3870719001172923
$ws.Range("M9").Value = 101.5800373333333
$ws.Range("N9").Value = 304.740112
$ws.Range("O9").Value = 0.3167483425780597
$ws.Range("P9").Value = 0.3167483425780597
$ws.Range("Q9").Value = 16098.18612266648
$ws.Range("R9").Value = 144883.6751039983
$ws.Range("S9").Value = 0.1226043828206926
$ws.Range("T9").Value = 0.1226043828206926

$ws.Range("G10").Value = 158.477852
$ws.Range("H10").Value = 475.433556
$ws.Range("I10").Value = 0.3870719001172923
$ws.Range("J10").Value = 0.3870719001172923
$ws.Range("M10").Value = 142.0267893333333
$ws.Range("N10").Value = 426.080368
$ws.Range("O10").Value = 0.4428699900492579
$ws.Range("P10").Value = 0.4428699900492579
$ws.Range("Q10").Value = 22508.10050000318
$ws.Range("R10").Value = 202572.9045000286
$ws.Range("S10").Value = 0.1714225285532926
$ws.Range("T10").Value = 0.1714225285532926
